# Generate Report for Handback
# Adds a new row (for file b4cd2c1f-232f-4643-b458-edbe832cb248.md) to each of
# the three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId = "b4cd2c1f-232f-4643-b458-edbe832cb248"
$mdName = "$fileId.md"
$mdPath = "e2e\$fileId.md"
$ghUrlMd = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$fileId.md"

$statusInSync = "Handed back: in sync with en-US"

$zhXlf = "$fileId.6827eed636b21d7975f2ae4177dc3bc2c5f1985f.zh-cn.xlf"
$deXlf = "$fileId.6827eed636b21d7975f2ae4177dc3bc2c5f1985f.de-de.xlf"

$zhHandoffDate = "2016-11-14 17:53:03"
$zhHandbackDate = "2016-11-14 17:53:49"
$deHandoffDate = "2016-11-14 17:53:17"
$deHandbackDate = "2016-11-14 17:54:07"
$overviewDate = "2016-11-14 17:53:17"

$zhGhUrlXlf = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000000/e2e/$fileId.md"
$deGhUrlXlf = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000000/e2e/$fileId.md"

# ---------------------------------------------------------------------------
# Sheet: Overview (table "Overview" -> columns A:G)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()
$rOverview = $rowOverview.Range.Row

$wsOverview.Cells.Item($rOverview, 1).Value = $mdName
$wsOverview.Cells.Item($rOverview, 1).Style = "Normal"

$wsOverview.Cells.Item($rOverview, 2).Value = $mdPath
$h1 = $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rOverview, 2), $ghUrlMd, "", "", $mdPath)

$wsOverview.Cells.Item($rOverview, 3).Value = ".md"
$wsOverview.Cells.Item($rOverview, 5).Value = $statusInSync
$wsOverview.Cells.Item($rOverview, 6).Value = $statusInSync
$wsOverview.Cells.Item($rOverview, 7).Value = $overviewDate
$wsOverview.Cells.Item($rOverview, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet: zh-cn (table "zh-cn" -> columns A:P)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")
$rowZh = $loZh.ListRows.Add()
$rZh = $rowZh.Range.Row

$wsZh.Cells.Item($rZh, 1).Value = $mdName
$h2 = $wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 1), $ghUrlMd, "", "", $mdName)

$wsZh.Cells.Item($rZh, 2).Value = ".md"
$wsZh.Cells.Item($rZh, 3).Value = $statusInSync
$wsZh.Cells.Item($rZh, 4).Value = "e2e"
$wsZh.Cells.Item($rZh, 5).Value = "ht"
$wsZh.Cells.Item($rZh, 6).Value = "True"
$wsZh.Cells.Item($rZh, 7).Value = $zhXlf
$wsZh.Cells.Item($rZh, 8).Value = $zhHandoffDate
$wsZh.Cells.Item($rZh, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Cells.Item($rZh, 9).Value = $mdName
$h3 = $wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 9), $zhGhUrlXlf, "", "", $mdName)

$wsZh.Cells.Item($rZh, 10).Value = $zhXlf
$wsZh.Cells.Item($rZh, 11).Value = $zhHandbackDate
$wsZh.Cells.Item($rZh, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($rZh, 12).Value = ""
$wsZh.Cells.Item($rZh, 13).Value = "True"
$wsZh.Cells.Item($rZh, 14).Value = ""
$wsZh.Cells.Item($rZh, 15).Value = "False"
$wsZh.Cells.Item($rZh, 16).Value = ""

# ---------------------------------------------------------------------------
# Sheet: de-de (table "de-de" -> columns A:P)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")
$rowDe = $loDe.ListRows.Add()
$rDe = $rowDe.Range.Row

$wsDe.Cells.Item($rDe, 1).Value = $mdName
$h4 = $wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 1), $ghUrlMd, "", "", $mdName)

$wsDe.Cells.Item($rDe, 2).Value = ".md"
$wsDe.Cells.Item($rDe, 3).Value = $statusInSync
$wsDe.Cells.Item($rDe, 4).Value = "e2e"
$wsDe.Cells.Item($rDe, 5).Value = "ht"
$wsDe.Cells.Item($rDe, 6).Value = "True"
$wsDe.Cells.Item($rDe, 7).Value = $deXlf
$wsDe.Cells.Item($rDe, 8).Value = $deHandoffDate
$wsDe.Cells.Item($rDe, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Cells.Item($rDe, 9).Value = $mdName
$h5 = $wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 9), $deGhUrlXlf, "", "", $mdName)

$wsDe.Cells.Item($rDe, 10).Value = $deXlf
$wsDe.Cells.Item($rDe, 11).Value = $deHandbackDate
$wsDe.Cells.Item($rDe, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($rDe, 12).Value = ""
$wsDe.Cells.Item($rDe, 13).Value = "True"
$wsDe.Cells.Item($rDe, 14).Value = ""
$wsDe.Cells.Item($rDe, 15).Value = "False"
$wsDe.Cells.Item($rDe, 16).Value = ""

Write-Host "Done adding handback rows."
